$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry: 2023-01-31 (serial 44957)
$ws.Range("A6").Value = 44957
$ws.Range("A6").NumberFormat = "d-mmm"
$ws.Range("B6").Value = "14.45 - 16.45"
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = "Continue researching algorithms. Implement cell neighbour finding"

# New entry: 2023-02-01 (serial 44958)
$ws.Range("A7").Value = 44958
$ws.Range("A7").NumberFormat = "d-mmm"
$ws.Range("B7").Value = "11.30 - 15.00"
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = "Implementing dijkstra algorithm"

# Update the active selection to reflect where the user was working
$ws.Range("E8").Select()
